$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("planned_operation_time", "str"),
    @("requested_service_url", "str"),
    @("time:timestamp", "datetime"),
    @("SubProcessID", "str"),
    @("unsatisfied_condition_description", "str"),
    @("parameters", "dict"),
    @("lifecycle:transition", "str"),
    @("org:resource", "str"),
    @("lifecycle:state", "str"),
    @("case:concept:name", "str"),
    @("operation_end_time", "datetime"),
    @("complete_service_time", "str"),
    @("current_task", "str"),
    @("human_workstation_green_button_pressed", "float"),
    @("event_id", "str"),
    @("response_status_code", "float"),
    @("concept:name", "str"),
    @("case", "str"),
    @("identifier:id", "str"),
    @("process_model_id", "str")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
